$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 827.0714
$ws.Range("I18").Value = 798.0909
$ws.Range("J18").Value = 933.3333
$ws.Range("K18").Value = 798.0909
$ws.Range("L18").Value = 933.3333
$ws.Range("M18").Value = -514.0909
$ws.Range("N18").Value = -1501.3333

$ws.Range("H21").Value = 17500
$ws.Range("I21").Value = 17500
$ws.Range("K21").Value = 17500
$ws.Range("M21").Value = -17032

$ws.Range("H23").Value = 17500
$ws.Range("I23").Value = 17500
$ws.Range("K23").Value = 17500
$ws.Range("M23").Value = -17266

$ws.Range("H29").Value = 650.55554
$ws.Range("I29").Value = 361
$ws.Range("K29").Value = 1083
$ws.Range("M29").Value = -802

$ws.Range("H33").Value = 760.05554
$ws.Range("I33").Value = 1053.7273
$ws.Range("J33").Value = 298.57144
$ws.Range("K33").Value = 1053.7273
$ws.Range("L33").Value = 298.57144
$ws.Range("M33").Value = -824.7273
$ws.Range("N33").Value = -756.5714399999999

$ws.Range("H38").Value = 3744.3125
$ws.Range("I38").Value = 2500
$ws.Range("J38").Value = 3922.0715
$ws.Range("K38").Value = 7500
$ws.Range("L38").Value = 11766.2145
$ws.Range("M38").Value = -7128
$ws.Range("N38").Value = -12510.2145

$ws.Range("H87").Value = 41232.875
$ws.Range("J87").Value = 41232.875
$ws.Range("L87").Value = 41232.875
$ws.Range("N87").Value = -43728.875

$ws.Range("H90").Value = 41232.875
$ws.Range("J90").Value = 41232.875
$ws.Range("L90").Value = 123698.625
$ws.Range("N90").Value = -136178.625

$ws.Range("H116").Value = 4797.0557
$ws.Range("I116").Value = 3017.7778
$ws.Range("J116").Value = 6576.3335
$ws.Range("K116").Value = 3017.7778
$ws.Range("L116").Value = 6576.3335
$ws.Range("M116").Value = 424.2222000000002
$ws.Range("N116").Value = -13460.3335

$ws.Range("H141").Value = 2262.4473
$ws.Range("I141").Value = 1732.2174
$ws.Range("K141").Value = 5196.6522
$ws.Range("M141").Value = -16.65220000000045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3476.6667
$ws.Range("I32").Value = 3306.1458
$ws.Range("J32").Value = 8933.333000000001
$ws.Range("K32").Value = 3306.1458
$ws.Range("L32").Value = 8933.333000000001
$ws.Range("M32").Value = -3019.1458
$ws.Range("N32").Value = -9507.333000000001

$ws.Range("H37").Value = 22351.188
$ws.Range("J37").Value = 23401.357
$ws.Range("L37").Value = 23401.357
$ws.Range("N37").Value = -23947.357

$ws.Range("H55").Value = 21744
$ws.Range("J55").Value = 39640
$ws.Range("L55").Value = 39640
$ws.Range("N55").Value = -40270

$ws.Range("H63").Value = 4800
$ws.Range("I63").Value = 2750
$ws.Range("J63").Value = 5971.4287
$ws.Range("K63").Value = 2750
$ws.Range("L63").Value = 5971.4287
$ws.Range("M63").Value = -2064
$ws.Range("N63").Value = -7343.4287

$ws.Range("H66").Value = 4800
$ws.Range("I66").Value = 2750
$ws.Range("J66").Value = 5971.4287
$ws.Range("K66").Value = 13750
$ws.Range("L66").Value = 29857.1435
$ws.Range("M66").Value = -10318
$ws.Range("N66").Value = -36721.14350000001

$ws.Range("H80").Value = 31450
$ws.Range("J80").Value = 38900
$ws.Range("L80").Value = 38900
$ws.Range("N80").Value = -40896

$ws.Range("H83").Value = 31450
$ws.Range("J83").Value = 38900
$ws.Range("L83").Value = 116700
$ws.Range("N83").Value = -126684

$ws.Range("H122").Value = 12635.895
$ws.Range("I122").Value = 14966.8
$ws.Range("J122").Value = 3895
$ws.Range("K122").Value = 44900.39999999999
$ws.Range("L122").Value = 11685
$ws.Range("M122").Value = -42450.39999999999
$ws.Range("N122").Value = -16585

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 25937.223
$ws.Range("J82").Value = 35616
$ws.Range("L82").Value = 35616
$ws.Range("N82").Value = -36382

$ws.Range("H85").Value = 25937.223
$ws.Range("J85").Value = 35616
$ws.Range("L85").Value = 35616
$ws.Range("N85").Value = -38268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9300.666999999999
$ws.Range("J50").Value = 9300.666999999999
$ws.Range("L50").Value = 9300.666999999999
$ws.Range("N50").Value = -10550.667

$ws.Range("H51").Value = 9256.571
$ws.Range("J51").Value = 9256.571
$ws.Range("L51").Value = 9256.571
$ws.Range("N51").Value = -10728.571

$ws.Range("H59").Value = 16295.857
$ws.Range("J59").Value = 16295.857
$ws.Range("L59").Value = 16295.857
$ws.Range("N59").Value = -18585.857

$ws.Range("H60").Value = 24543.5
$ws.Range("J60").Value = 24543.5
$ws.Range("L60").Value = 24543.5
$ws.Range("N60").Value = -25565.5

$ws.Range("H61").Value = 9256.571
$ws.Range("J61").Value = 9256.571
$ws.Range("L61").Value = 9256.571
$ws.Range("N61").Value = -9952.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1261.7273
$ws.Range("I34").Value = 219.16667
$ws.Range("J34").Value = 1652.6875
$ws.Range("K34").Value = 657.50001
$ws.Range("L34").Value = 4958.0625
$ws.Range("M34").Value = -573.50001
$ws.Range("N34").Value = -5126.0625

$ws.Range("H39").Value = 2327.7
$ws.Range("J39").Value = 2327.7
$ws.Range("L39").Value = 6983.099999999999
$ws.Range("N39").Value = -7571.099999999999

$ws.Range("H113").Value = 363095.22
$ws.Range("I113").Value = 611.3889
$ws.Range("J113").Value = 864995.9399999999
$ws.Range("K113").Value = 1834.1667
$ws.Range("L113").Value = 2594987.82
$ws.Range("M113").Value = 335.8332999999998
$ws.Range("N113").Value = -2599327.82

$ws.Range("H130").Value = 2209.476
$ws.Range("I130").Value = 879.8
$ws.Range("J130").Value = 2625
$ws.Range("K130").Value = 2639.4
$ws.Range("L130").Value = 7875
$ws.Range("M130").Value = 2380.6
$ws.Range("N130").Value = -17915

$ws.Range("H131").Value = 892.92
$ws.Range("I131").Value = 478.46155
$ws.Range("J131").Value = 954.8506
$ws.Range("K131").Value = 1435.38465
$ws.Range("L131").Value = 2864.5518
$ws.Range("M131").Value = 3604.61535
$ws.Range("N131").Value = -12944.5518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 829.3333
$ws.Range("I22").Value = 829.3333
$ws.Range("K22").Value = 829.3333
$ws.Range("M22").Value = -300.3333

$ws.Range("H43").Value = 4486.5
$ws.Range("J43").Value = 9507.6
$ws.Range("L43").Value = 9507.6
$ws.Range("N43").Value = -9809.6

$ws.Range("H46").Value = 10272.533
$ws.Range("J46").Value = 12228.8
$ws.Range("L46").Value = 12228.8
$ws.Range("N46").Value = -12540.8

$ws.Range("H57").Value = 14174.637
$ws.Range("J57").Value = 17663.875
$ws.Range("L57").Value = 17663.875
$ws.Range("N57").Value = -19303.875

$ws.Range("H70").Value = 41079.465
$ws.Range("I70").Value = 48996.348
$ws.Range("K70").Value = 48996.348
$ws.Range("M70").Value = -48726.348

$ws.Range("H73").Value = 41079.465
$ws.Range("I73").Value = 48996.348
$ws.Range("K73").Value = 48996.348
$ws.Range("M73").Value = -48060.348

$ws.Range("H80").Value = 2671.4285
$ws.Range("I80").Value = 2350
$ws.Range("K80").Value = 2350
$ws.Range("M80").Value = -1352

$ws.Range("H83").Value = 2671.4285
$ws.Range("I83").Value = 2350
$ws.Range("K83").Value = 11750
$ws.Range("M83").Value = -6758

$ws.Range("H122").Value = 2175.4187
$ws.Range("I122").Value = 2083.9656
$ws.Range("J122").Value = 2364.8572
$ws.Range("K122").Value = 6251.8968
$ws.Range("L122").Value = 7094.571599999999
$ws.Range("M122").Value = -3801.8968
$ws.Range("N122").Value = -11994.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2200
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 2644.4443
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 2644.4443
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -3020.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9466.375
$ws.Range("J54").Value = 9466.375
$ws.Range("L54").Value = 9466.375
$ws.Range("N54").Value = -10506.375

$ws.Range("H81").Value = 8425.532999999999
$ws.Range("J81").Value = 1996
$ws.Range("L81").Value = 3992
$ws.Range("N81").Value = -6114

$ws.Range("H84").Value = 8425.532999999999
$ws.Range("J84").Value = 1996
$ws.Range("L84").Value = 19960
$ws.Range("N84").Value = -30568

$ws.Range("H113").Value = 56944.832
$ws.Range("I113").Value = 67066.87
$ws.Range("J113").Value = 6334.6665
$ws.Range("K113").Value = 201200.61
$ws.Range("L113").Value = 19003.9995
$ws.Range("M113").Value = -199030.61
$ws.Range("N113").Value = -23343.9995

$ws.Range("H132").Value = 1848.6522
$ws.Range("I132").Value = 1656.7115
$ws.Range("J132").Value = 2435.7646
$ws.Range("K132").Value = 4970.1345
$ws.Range("L132").Value = 7307.293799999999
$ws.Range("M132").Value = -2440.1345
$ws.Range("N132").Value = -12367.2938
